# Daily attendance processing - 2026-01-15 10:39:10
#
# For every row in the "Recorded By" column (G), when the recorded-by
# list starts with "System" (or "system"), cycle that leading name to
# the end of the comma-separated list, e.g.
#   "System, dnasr281@gmail.com"                -> "dnasr281@gmail.com, System"
#   "system, backup@backdoor.com, System"       -> "backup@backdoor.com, System, system"
# Rows whose list does not start with System/system (single names, or
# lists starting with another name) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count
$col = 7  # column G - "Recorded By"

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $col)
    $val = $cell.Text

    if ($val -eq $null -or $val -eq "") {
        continue
    }

    $parts = $val.Split(",")
    if ($parts.Length -le 1) {
        continue
    }

    $first = $parts[0].Trim()
    if ($first -ne "System" -and $first -ne "system") {
        continue
    }

    $newParts = @()
    for ($i = 1; $i -lt $parts.Length; $i++) {
        $newParts += $parts[$i].Trim()
    }
    $newParts += $first

    $cell.Value = ($newParts -join ", ")
}
